$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that ends the abstract body text ("The contributions
# of this thesis represent advancements ...") - the new "Contributions" list
# is appended right after it, before the pre-existing trailing blank
# paragraphs.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*The contributions of this thesis represent advancements*") {
        $target = $p
    }
}

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Creates a brand-new empty paragraph right after $para and returns it.
# New paragraphs otherwise inherit $para's explicit paragraph style (instead
# of resolving the style's "next" style the way interactive typing would),
# so it is reset back to Normal/Default here; callers that want a styled
# paragraph set the style explicitly afterwards anyway.
function New-ParaAfter($para) {
    $r = $para.Range.Duplicate
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newIndex = $para.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Style = "Normal"
    return $newPara
}

# Returns a zero-length Range positioned just before the paragraph's own
# trailing paragraph mark - safe to call InsertAfter/InsertBreak on without
# the text/break leaking into the following paragraph (collapsing exactly on
# a paragraph boundary is otherwise "sticky" toward the next paragraph).
function Get-EndInsertionPoint($para) {
    $e = $para.Range.End
    return $d.Range($e - 1, $e - 1)
}

# Sets the paragraph's style and appends text into it.
function Set-ParaStyleAndText($para, $styleName, $text) {
    $para.Style = $styleName
    $ip = Get-EndInsertionPoint $para
    $ip.InsertAfter($text)
}

# ---------------------------------------------------------------------------
# Build the new content, paragraph by paragraph, exactly mirroring the diff:
#
#   <blank>
#   <blank>
#   <left-aligned paragraph containing a page break>
#   Heading3: "Contributions:"
#   <blank>
#   Heading4: " ArXiV: Radiotherapy Dosimetry: A Review on Open-Source Optimizer"
#   <blank>
#   Heading4: "ESTRO: A Novel Framework for Multi-Objective Optimization and
#              Robust Plan Selection Using Graph Theory"
#   <blank>
#   Heading4: "SFPM: Dose Volume Histograms Guided Deep Dose Predictions"
#   <blank>
#   Heading4: "AIME: Radiotherapy Dose Optimization via Clinical Knowledge
#              Based Reinforcement Learning (full paper coming soon)"
#   <blank>
#   Heading4: "ASTRO: Clinically Dependent Fully Automatic Treatment Planning
#              System"
#   <blank>
#   Heading4: "SFRO: Attention Mechanism on Dose-Volume Histograms for Deep
#              Dose Predictions"
#
# Note: InsertBreak(wdPageBreak) both places the break character inside the
# *current* paragraph and implicitly starts a brand-new paragraph after it
# (mirroring a manual Ctrl+Enter keypress), so the "blank2" paragraph is
# reused as the page-break holder instead of adding yet another paragraph.
# ---------------------------------------------------------------------------

$pBlank1 = New-ParaAfter $target
$pBlank2 = New-ParaAfter $pBlank1

$breakIp = Get-EndInsertionPoint $pBlank2
$breakIp.InsertBreak(7)
$pBreak = $d.Paragraphs.Item($pBlank2.Index)
$pBreak.Alignment = 0

$pContrib = $d.Paragraphs.Item($pBreak.Index + 1)
Set-ParaStyleAndText $pContrib "Heading 3" "Contributions:"

$pBlank3 = New-ParaAfter $pContrib

$pArxiv = New-ParaAfter $pBlank3
Set-ParaStyleAndText $pArxiv "Heading 4" " ArXiV: Radiotherapy Dosimetry: A Review on Open-Source Optimizer"

$pBlank4 = New-ParaAfter $pArxiv

$pEstro = New-ParaAfter $pBlank4
Set-ParaStyleAndText $pEstro "Heading 4" "ESTRO: A Novel Framework for Multi-Objective Optimization and Robust Plan Selection Using Graph Theory"

$pBlank5 = New-ParaAfter $pEstro

$pSfpm = New-ParaAfter $pBlank5
Set-ParaStyleAndText $pSfpm "Heading 4" "SFPM: Dose Volume Histograms Guided Deep Dose Predictions"

$pBlank6 = New-ParaAfter $pSfpm

$pAime = New-ParaAfter $pBlank6
Set-ParaStyleAndText $pAime "Heading 4" "AIME: Radiotherapy Dose Optimization via Clinical Knowledge Based Reinforcement Learning (full paper coming soon)"

$pBlank7 = New-ParaAfter $pAime

$pAstro = New-ParaAfter $pBlank7
Set-ParaStyleAndText $pAstro "Heading 4" "ASTRO: Clinically Dependent Fully Automatic Treatment Planning System"

$pBlank8 = New-ParaAfter $pAstro

$pSfro = New-ParaAfter $pBlank8
Set-ParaStyleAndText $pSfro "Heading 4" "SFRO: Attention Mechanism on Dose-Volume Histograms for Deep Dose Predictions"
